$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Planilha1")

# Clear old contents for rows 3-6 across columns A:C, then set new values
$ws.Range("A3:C6").ClearContents()

$ws.Range("A3").Value = "Esquadria fora de prumo"
$ws.Range("B3").Value = "Esquadrias"

$ws.Range("A4").Value = "Instalações aparentes"
$ws.Range("B4").Value = "Elétrica"

$ws.Range("A5").Value = "Pintura irregular"
$ws.Range("B5").Value = "Conservação"

$ws.Range("A6").Value = "Resíduos"
$ws.Range("B6").Value = "Conservação"
